$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns for rows with refreshed crypto data
$ws.Range("D2").Value = "26.062.76"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.650.81"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5261"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06322"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07799"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.507"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "1.672.64"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5487"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "0.0₅8211"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "26.076.64"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.572"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.038"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1237"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.226"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05807"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.273"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.540"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.258"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.585"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.411"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9447"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.774"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01607"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8439"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.744"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.11%  "
$ws.Range("D43").Value = "1.027.78"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").Value = "1.795.99"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4313"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09641"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.60%  "

# Rows 48-50 were reordered: EnergySwap, Cronos, RenderToken now rank 48-50
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.857"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05148"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.467"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.27%  "
